$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows per repull/mean recalculation
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -7
$ws.Range("F6").Value = -5
$ws.Range("F8").Value = -4
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = -5
$ws.Range("F13").Value = -8
